# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap country names in A206 and A207 (Timor Oriental <-> Santa Lucia)
$a206 = $ws.Range("A206").Value2
$a207 = $ws.Range("A207").Value2
$ws.Range("A206").Value2 = $a207
$ws.Range("A207").Value2 = $a206

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 20:59"

# Update numeric data for the changed country rows
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7268282
$ws.Range("C4").Value = 24098
$ws.Range("D4").Value = 4502646
$ws.Range("E4").Value = 2556801
$ws.Range("G4").Value = 395
$ws.Range("H4").Value = 208835

# Row 25 - Alemania
$ws.Range("B25").Value = 284984
$ws.Range("C25").Value = 1278
$ws.Range("E25").Value = 25952
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9532

# Row 27 - Israel
$ws.Range("B27").Value = 226586
$ws.Range("C27").Value = 8687
$ws.Range("D27").Value = 157523
$ws.Range("E27").Value = 67646
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 1417

# Row 36 - Republica Dominicana
$ws.Range("B36").Value = 110957
$ws.Range("C36").Value = 360
$ws.Range("D36").Value = 85220
$ws.Range("E36").Value = 23644
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 2093

# Row 56 - (no country name change)
$ws.Range("E56").Value = 6476
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 239

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 54819
$ws.Range("C60").Value = 427
$ws.Range("D60").Value = 51322
$ws.Range("E60").Value = 3045
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 452

# Row 76 - Libia
$ws.Range("B76").Value = 31828
$ws.Range("C76").Value = 538
$ws.Range("E76").Value = 13821
$ws.Range("G76").Value = 8
$ws.Range("H76").Value = 499

# Row 112 - Mozambique
$ws.Range("B112").Value = 7757
$ws.Range("C112").Value = 168
$ws.Range("D112").Value = 4769
$ws.Range("E112").Value = 2934
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 54

# Row 135 - Aruba
$ws.Range("B135").Value = 3832
$ws.Range("C135").Value = 33
$ws.Range("D135").Value = 2829
$ws.Range("E135").Value = 978

# Row 144 - Mali
$ws.Range("B144").Value = 3080
$ws.Range("C144").Value = 16
$ws.Range("D144").Value = 2410
$ws.Range("E144").Value = 540
